# Generate Report for Handback
# Updates timestamps / priority values produced by a re-run of the
# handback report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" for the
#     7ce0f49b-... row (row 3) moved from 00:24:03 to 00:24:58.
$wsOverview.Range("G3").Value = "2016-09-02 00:24:58"
$wsOverview.Range("G4").Value = "2016-09-02 00:24:58"

# --- zh-cn sheet, row for 7ce0f49b-... (row 3):
#   Priority: ht -> mt
#   Correspond Handoff Datetime:  2016-09-02 00:23:56 -> 2016-09-02 00:24:54
#   Correspond Handback DateTime: 2016-09-02 00:24:28 -> 2016-09-02 00:25:18
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-02 00:24:54"
$wsZhCn.Range("H4").Value = "2016-09-02 00:24:54"
$wsZhCn.Range("K3").Value = "2016-09-02 00:25:18"
$wsZhCn.Range("K4").Value = "2016-09-02 00:25:18"

# --- de-de sheet, row for 7ce0f49b-... (row 3):
#   Priority: ht -> mt
#   Correspond Handoff Datetime:  2016-09-02 00:24:03 -> 2016-09-02 00:24:58
#   Correspond Handback DateTime: 2016-09-02 00:24:35 -> 2016-09-02 00:25:25
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-02 00:24:58"
$wsDeDe.Range("H4").Value = "2016-09-02 00:24:58"
$wsDeDe.Range("K3").Value = "2016-09-02 00:25:25"
$wsDeDe.Range("K4").Value = "2016-09-02 00:25:25"
